# Update cryptos list price / volume(1h) figures.
# Note: several "Price" (column D) values are plain-looking decimal numbers
# (e.g. "23.85"). Excel's COM layer auto-converts such text into a numeric
# value on assignment, which would change the cell's stored type away from
# text. Since the source data keeps these as text strings, we force the
# cell's number format to Text ("@") immediately before writing any such
# value so it is stored as a string, matching the original file's cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.209.52"
$ws.Range("E2").Value = "  +2.49%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.586.94"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +1.22%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.22"
$ws.Range("E5").Value = "  +1.31%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.45%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +1.21%  "

# Row 8 - Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.85"
$ws.Range("E8").Value = "  +6.13%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.13%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.29%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("E11").Value = "  +2.33%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.814.51"
$ws.Range("E12").Value = "  +1.50%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.584.48"
$ws.Range("E13").Value = "  +0.89%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  +1.42%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -0.22%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "28.230.29"
$ws.Range("E16").Value = "  +2.62%  "

# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.11"
$ws.Range("E17").Value = "  +1.08%  "

# Row 18 - BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.13"
$ws.Range("E18").Value = "  +1.13%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +0.25%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -0.91%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +1.25%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("E22").Value = "  -1.62%  "

# Row 23 - Avalanche
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("E23").Value = "  -1.17%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.56%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.93"
$ws.Range("E25").Value = "  +1.10%  "

# Row 26 - EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.13"
$ws.Range("E26").Value = "  -0.25%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.60%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  -1.04%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +1.23%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.29%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("E31").Value = "  +0.00%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.31%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -1.27%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.398.89"
$ws.Range("E34").Value = "  -4.34%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  -7.52%  "

# Row 37 - HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  +1.54%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.41%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  +8.85%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  -0.31%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  -0.87%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +1.21%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  -0.31%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  -1.91%  "

# Row 45 - WEMIXToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.980"
$ws.Range("E45").Value = "  +0.62%  "

# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.15"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.724.17"
$ws.Range("E47").Value = "  +1.24%  "

# Row 48 - mCoin
$ws.Range("E48").Value = "  +2.02%  "

# Row 49 - Quant
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.89"
$ws.Range("E49").Value = "  +0.35%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +5.95%  "

# Row 51 - Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0520"
$ws.Range("E51").Value = "  -0.69%  "
